# Updated cryptos list (scraped price/volume refresh).
# For "Price" (column D) cells whose new value parses as a plain decimal
# number (e.g. "7.48", "1.00"), NumberFormat is forced to Text ("@") right
# before the assignment so Excel stores the literal digits/trailing zeros
# instead of silently coercing the cell to a numeric value. Multi-dot
# values (e.g. "70.764.44") and the percent strings in column E are never
# numeric-looking, so they are assigned directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.764.44'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').Value = '3.809.09'
$ws.Range('E3').Value = '  -1.15%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '709.74'
$ws.Range('E5').Value = '  +1.80%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '170.61'
$ws.Range('E6').Value = '  -1.48%  '
$ws.Range('D7').Value = '3.809.23'
$ws.Range('E7').Value = '  -1.11%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -0.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.160'
$ws.Range('E10').Value = '  -1.59%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.48'
$ws.Range('E11').Value = '  +4.96%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.458'
$ws.Range('E12').Value = '  -0.49%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000253'
$ws.Range('E13').Value = '  -2.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.00'
$ws.Range('E14').Value = '  -1.01%  '
$ws.Range('D15').Value = '4.453.84'
$ws.Range('E15').Value = '  -1.16%  '
$ws.Range('D16').Value = '3.832.94'
$ws.Range('E16').Value = '  -0.52%  '
$ws.Range('D17').Value = '70.850.28'
$ws.Range('E17').Value = '  -0.53%  '
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('E19').Value = '  -1.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.31'
$ws.Range('E20').Value = '  -2.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '494.90'
$ws.Range('E21').Value = '  +0.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.65'
$ws.Range('E22').Value = '  -4.41%  '
$ws.Range('E23').Value = '  +0.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.33'
$ws.Range('E24').Value = '  -0.98%  '
$ws.Range('E25').Value = '  -1.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.08'
$ws.Range('E26').Value = '  -1.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.41'
$ws.Range('E27').Value = '  -2.45%  '
$ws.Range('D28').Value = '3.962.62'
$ws.Range('E28').Value = '  -1.24%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('E30').Value = '  -4.81%  '
$ws.Range('E31').Value = '  -3.00%  '
$ws.Range('E32').Value = '  -1.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.34'
$ws.Range('E33').Value = '  -4.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.09'
$ws.Range('E34').Value = '  -1.95%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.172'
$ws.Range('E35').Value = '  -3.18%  '
$ws.Range('B36').Value = 'RenzoRestakedETH'
$ws.Range('C36').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D36').Value = '3.779.65'
$ws.Range('E36').Value = '  -0.70%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '9.11'
$ws.Range('E37').Value = '  -1.78%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('E39').Value = '  -2.47%  '
$ws.Range('E40').Value = '  +0.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.30'
$ws.Range('E41').Value = '  -3.46%  '
$ws.Range('E42').Value = '  -1.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.24'
$ws.Range('E43').Value = '  -4.50%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  +0.10%  '
$ws.Range('B46').Value = 'FLOKI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.000321'
$ws.Range('E46').Value = '  +3.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '166.16'
$ws.Range('E47').Value = '  +1.60%  '
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '425.68'
$ws.Range('E48').Value = '  +1.43%  '
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '48.81'
$ws.Range('E49').Value = '  +0.32%  '
$ws.Range('E50').Value = '  -0.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.294'
$ws.Range('E51').Value = '  -2.82%  '
